# Applies the edit described in the commit:
# "Atualizacao aula 09 (que foi movida para 16/09)"
#
# Shifts the lesson-plan content of rows 7-14 down (the "Aula estudio
# para atividade Projeto Profissional" double lesson now also appears on
# 02/09, pushing Distribuicao/Internacionalizacao/Licencas down a slot and
# moving the "Aula estudio Testes de Software" studio lesson from 16/09 to
# 23/09), and re-applies identical formatting on the FERIADO rows (18, 24).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$FILL_YELLOW = 13957622   # FFF6F9D4 (light yellow-green highlight)
$FILL_RED = 4260031       # FFBF0041 (FERIADO dark red highlight)
$NO_FILL = $null

function Set-Cell {
    param($addr, $text, $wrap, $fillColor)
    $r = $ws.Range($addr)
    $r.Value = $text
    $r.WrapText = $wrap
    if ($fillColor -eq $null) {
        $r.Interior.ColorIndex = -4142
    } else {
        $r.Interior.Color = $fillColor
    }
}

# --- Row heights that change ---
$ws.Rows.Item(7).RowHeight = 89.25
$ws.Rows.Item(8).RowHeight = 13.300000000000001
$ws.Rows.Item(9).RowHeight = 76.5
$ws.Rows.Item(10).RowHeight = 165.75
$ws.Rows.Item(11).RowHeight = 63.75
$ws.Rows.Item(12).RowHeight = 63.75
$ws.Rows.Item(13).RowHeight = 51
$ws.Rows.Item(14).RowHeight = 51

# --- Cell content / formatting ---

# Row 7
Set-Cell 'B7' 'Aula estúdio para atividade “Projeto Profissional” ' $true $NO_FILL
Set-Cell 'C7' 'Aula estúdio para atividade “Projeto Profissional” ' $true $NO_FILL
Set-Cell 'D7' 'Atividades entregues via PR no github.' $true $NO_FILL
Set-Cell 'E7' 'Aula estúdio para atividade “Projeto Profissional” ' $true $NO_FILL

# Row 8
Set-Cell 'B8' 'FERIADO' $false $FILL_RED
Set-Cell 'C8' 'FERIADO' $false $FILL_RED
Set-Cell 'D8' 'FERIADO' $false $FILL_RED
Set-Cell 'E8' 'FERIADO' $false $FILL_RED

# Row 9
Set-Cell 'B9' 'Distribuição e versionamento de software' $true $FILL_YELLOW
Set-Cell 'C9' 'Ferramentas de distribuição de software. Tipos de distribuição. Esquemas de versionamento. Referência: FOGEL Cap 7' $true $NO_FILL
Set-Cell 'D9' 'Atividades entregues via PR no github.' $true $NO_FILL
Set-Cell 'E9' 'Expositiva; Atividade prática: criação de um pacote Python para distribuição via pip e de um pacote Debian para instalação em sistemas baseados em Debian.' $true $NO_FILL

# Row 10
Set-Cell 'B10' 'Como fazer software funcionar em vários países?' $true $FILL_YELLOW
Set-Cell 'C10' 'Definição de localização e internacionalização. Ferramenta gettext para traduções e format .po.' $true $NO_FILL
Set-Cell 'D10' 'Atividades entregues via PR no github.' $true $NO_FILL
Set-Cell 'E10' 'Expositiva; Atividade prática: tradução e localização de um aplicativo linha de comando e de um sistema web em Flask.' $true $NO_FILL

# Row 11
Set-Cell 'B11' 'Licenças de software' $true $FILL_YELLOW
Set-Cell 'C11' 'Direitos autorais no Brasil, patentes e licenças de software;
Referência: LAURENT Cap 1 e 5, STALLMAN Cap 6' $true $NO_FILL
Set-Cell 'D11' 'Atividades entregues via PR no github.' $true $NO_FILL
Set-Cell 'E11' 'Expositiva; Atividade prática: discussões de casos relevantes recentes' $true $NO_FILL

# Row 12
Set-Cell 'B12' 'Testes de software' $false $NO_FILL
Set-Cell 'C12' 'Tipos de testes de software: unitarios, de integração e de interface. Ferramentas de testes: pytest, selenium' $false $NO_FILL
Set-Cell 'D12' 'Atividades entregues via PR no github.' $true $NO_FILL
Set-Cell 'E12' 'Expositiva: Ideias por trás de testes de software. 
Prática: adiciona testes de unidade e de interface ao projeto da atividade passada' $true $NO_FILL

# Row 13
Set-Cell 'B13' 'Aula estúdio para atividade “Testes de Software” ' $true $NO_FILL
Set-Cell 'C13' 'Aula estúdio para atividade “Testes de Software” ' $true $NO_FILL
Set-Cell 'D13' 'Atividades entregues via PR no github.' $true $NO_FILL
Set-Cell 'E13' 'Aula estúdio para atividade “Testes de Software” ' $true $NO_FILL

# Row 14
Set-Cell 'B14' 'Introdução a contribuições em projetos abertos - Atividade "Minha primeira contribuição"' $true $FILL_YELLOW

# Row 18
Set-Cell 'B18' 'FERIADO' $true $FILL_RED
Set-Cell 'C18' 'FERIADO' $true $FILL_RED
Set-Cell 'D18' 'FERIADO' $true $FILL_RED
Set-Cell 'E18' 'FERIADO' $true $FILL_RED

# Row 24
Set-Cell 'B24' 'FERIADO' $true $FILL_RED
Set-Cell 'C24' 'FERIADO' $true $FILL_RED
Set-Cell 'D24' 'FERIADO' $true $FILL_RED
Set-Cell 'E24' 'FERIADO' $true $FILL_RED
